# Daily attendance processing - 2025-12-10 15:56:20
#
# Normalizes the "Recorded By" column (G) so that multi-author entries are
# listed in reverse order (e.g. "a@b.com, System" -> "System, a@b.com").
# Single-author cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$firstRow = $used.Row

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val.Split(",")
        $n = $parts.Length

        $reversed = @()
        for ($j = $n - 1; $j -ge 0; $j--) {
            $reversed += $parts[$j].Trim()
        }

        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
